$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Customer has ( Id , Addresses , ..." -> insert "First Line ," after "Id ,"
# ---------------------------------------------------------------------------
$pCustomer = $d.Paragraphs.Item(2)
$r = $pCustomer.Range.Duplicate
$found = $r.Find.Execute("Customer has ( Id , Addresses , ", $true, $false, $false, $false, $false, $true, 1, $false, "Customer has ( Id ,First Line , Addresses , ", 2)

# ---------------------------------------------------------------------------
# 2. "Order details ( Id )" -> "Order details ( Id , Quantity , Unit Price )"
# ---------------------------------------------------------------------------
$pOrderDetails = $d.Paragraphs.Item(13)
$r = $pOrderDetails.Range.Duplicate
$found = $r.Find.Execute("Order details ( Id )", $true, $false, $false, $false, $false, $true, 1, $false, "Order details ( Id , Quantity , Unit Price )", 2)

# ---------------------------------------------------------------------------
# 3. Append new entity paragraphs (numId=2) after "Shipment ( Id, Tracking ... )"
# ---------------------------------------------------------------------------
$newEntities = @(
    "Country ( Id , Code ,Name )",
    "City ( Id, Code , Name )",
    "Zone ( Id , Code , Name )",
    "Distinct ( Id , Code , Name )",
    "Product Unit Details ( Id )",
    "Color ( Id , Name , Created At , Updated At )",
    "Size ( Id , Value , Created At , Updated At )"
)

$idx = 15
foreach ($txt in $newEntities) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $newP = $d.Paragraphs.Item($idx)
    $newP.Range.Text = $txt
}

# ---------------------------------------------------------------------------
# 4. Append new relationship paragraphs (numId=3) after
#    "Each Order must has one Shipment , Each Shipment must has one Order"
# ---------------------------------------------------------------------------
$relText = "Each Order must has one Shipment , Each Shipment must has one Order"
$relIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains($relText)) {
        $relIndex = $i
    }
}

$newRelations = @(
    "Each Customer may assigned to Many Distinct , Each Distinct may has many Customer",
    "Each country must has many Cities , Each City must has one Country",
    "Each City must has many Zones , Each Zone must has one City",
    "Each Zone must has many Distinct , each Distinct must has one zone",
    "Each Product Unit must has one Product Unit Details , each Product Unit Details must assign to Product Unit",
    "Each Color may has many Product Unit Details ,Product Unit Details may assigned to one color",
    "Each Size may has many Product Unit Details , Product Unit Details may assigned to one Size"
)

$idx = $relIndex
foreach ($txt in $newRelations) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $newP = $d.Paragraphs.Item($idx)
    $newP.Range.Text = $txt
}

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
